# RFWBBoard_RevB BOM update: add CJS-1200 SPDT switch line (SW1-SW4)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 36; everything from the old row 36 down shifts to row 37+.
# (This also auto-adjusts the SUM(J4:J37) total formula range.)
$ws.Rows(36).Insert()

# Populate the new row 36 with the switch BOM line.
$ws.Range("A36").Value = 33
$ws.Range("B36").Value = "SW1, SW2, SW3, SW4"
$ws.Range("C36").Value = "SWITCH-CJS-1200-1200"
$ws.Range("D36").Value = "J-Lead"
$ws.Range("E36").Value = "CJS"
$ws.Range("F36").Value = "CJS-1200TA"
$ws.Range("G36").Value = "Digi-Key"
$ws.Range("H36").Value = "563-1021-1-ND"
$ws.Range("I36").Value = "CJS-1200 SPDT Switch"
$ws.Range("J36").Value = 4

# Give the newly-inserted row (and the previously-formula-less rows above it)
# a Total Cost formula, like the rest of the BOM rows.
$ws.Range("L32:L36").Formula = "=J32*K32"

# The item numbers for the rows that shifted down need bumping by one.
$ws.Range("A37").Value = 34
$ws.Range("A38").Value = 35

# The old row 36's Total Cost formula shifted down to L37 along with the
# row; the final sheet no longer carries a formula there, so clear it.
$ws.Range("L37").ClearContents()
